# Rename the StockSex header labels (columns I1:V1) from the old
# "F"/"M"-prefixed naming convention to the new ".." subgroup naming
# convention (e.g. "FUPSALM" -> "UPSALM..F", "MUPSALM" -> "UPSALM..M").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "UPSALM..F"
$ws.Range("J1").Value = "UPSALM..M"
$ws.Range("K1").Value = "MFSALM..F"
$ws.Range("L1").Value = "MFSALM..M"
$ws.Range("M1").Value = "CHMBLN..F"
$ws.Range("N1").Value = "CHMBLN..M"
$ws.Range("O1").Value = "SFSALM..F"
$ws.Range("P1").Value = "SFSALM..M"
$ws.Range("Q1").Value = "HELLSC..F"
$ws.Range("R1").Value = "HELLSC..M"
$ws.Range("S1").Value = "TUCANO..F"
$ws.Range("T1").Value = "TUCANO..M"
$ws.Range("U1").Value = "FALL..F"
$ws.Range("V1").Value = "FALL..M"
